# TORIBDA feat: added calculateTotalPrice method
# Fill in row 6 (calculateTotalPrice) with its tracked time data,
# and move the active selection down to the next empty row (B7:E8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "5 mins"
$ws.Range("C6").Value = "1 min 28 seconds"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"

$ws.Range("B7:E8").Select()
